$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column before N (column 14). This shifts the old N,O,P,Q
#    (and anything further right) one column to the right, preserving their
#    exact widths/values/styles.
$ws.Columns("N").Insert()

# 2. Give the newly inserted column N the same width as column M (closest
#    achievable value in this engine's pixel-rounded ColumnWidth setter).
$ws.Columns("N").ColumnWidth = 16.667

# 3. The inserted column pushed the original "Message" column (header +
#    "Message for New transmittal" data) from N to O. Copy that content
#    back into N (duplicating it), so both N and O hold it for now.
$ws.Range("O1:O5").Copy()
$ws.Range("N1").PasteSpecial(-4104)

# 4. Column N now becomes the "<page> of <count> <message>" formula for the
#    data rows (row 1 keeps its original "Message" header untouched).
for ($r = 2; $r -le 5; $r++) {
    $ws.Range("N$r").Formula = "=CONCATENATE(ROW()-1,"" of "",COUNTA(A2:A100),"" "",O2)"
}

# 5. Update the _FilterDatabase defined name so it covers the new column R.
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=Transmittals_Close_Cancel!`$A`$1:`$R`$5"
    }
}
